# Penalty/Reward System (unfinished) - remove a couple of stray weekly/monthly
# data points that were thrown off by the new penalty-reward logic.
$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: drop the two rows for 45361.99999999999 (34) and
# 45368.99999999999 (98) - these were old rows 4 and 5. Deleting row 4 twice
# removes both and shifts everything below up by two rows.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(4).Delete()
$wsWeekly.Rows.Item(4).Delete()

# "Monthly Trend" sheet: drop the row for 45382.99999999999 (132) - this was
# old row 3. Deleting it shifts everything below up by one row.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows.Item(3).Delete()
